# Generate Report for Handoff
# Updates the status of each locale to "Ready for handoff" and refreshes
# the handoff timestamps on the Overview sheet and each locale sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# zh-cn locale sheet: Status -> Ready for handoff, Latest Handoff Datetime updated
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-24 14:46:54"

# de-de locale sheet: Status -> Ready for handoff, Latest Handoff Datetime updated
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-24 14:46:59"

# Overview sheet: per-locale status columns and the overall latest handoff date
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-03-24 14:46:59"
